# Fix the product-name value (typo: missing hyphen after the leading
# numeric code) on both the input and output sheets, and leave the
# workbook with the ProductLoanOutput sheet/cell selected - matching
# the state the workbook was saved in after the fix.

$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

$newProductName = "341-MS-EPP-DB-SAR-REC-NON-RNI-CTPD-SAR-MD-TR-1-ONTIME"

# Update the product name cell on the input sheet.
$wsInput.Activate()
$wsInput.Range("B1").Value = $newProductName
$wsInput.Range("B1").Select()

# Update the same product name cell on the output sheet.
$wsOutput.Activate()
$wsOutput.Range("B1").Value = $newProductName
$wsOutput.Range("B1").Select()

# ProductLoanOutput ends up as the active (selected) tab.
$wsOutput.Activate()
